# Update the public EPEX Spot prices workbook:
#  1. "Prix Spot": insert a new date column at EM (shifts 01-oct..31-oct right
#     by one column to EN..FR) and fill it with the new "10-dec" header plus
#     "-" placeholders for the 24 hourly rows (no data yet for that date).
#  2. "Gaz": append a new row for 2025-12-08 with its price.
#  3. "CO2": append a new row for 2025-12-08 with no price available yet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Prix Spot — insert column before EM, pushing the "01-oct." ... "31-oct."
#    block one column to the right.
# ---------------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")
$wsPrix.Range("EM1").EntireColumn.Insert()

$wsPrix.Range("EM1").Value = "10-dec"
$wsPrix.Range("EM1").Font.Bold = $true
$wsPrix.Range("EM1").HorizontalAlignment = -4108
$wsPrix.Range("EM1").VerticalAlignment = -4160
$wsPrix.Range("EM1").Borders.LineStyle = 1
$wsPrix.Range("EM1").Borders.Weight = 2

for ($row = 2; $row -le 25; $row++) {
    $wsPrix.Cells.Item($row, 143).Value = "-"
}

# ---------------------------------------------------------------------------
# 2) Gaz — append 2025-12-08 price row.
#    The date column is stored as literal text (not a real date) everywhere
#    else in the sheet, so force text formatting before assigning the value,
#    otherwise Excel auto-parses "2025-12-08" into a date serial.
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A173").NumberFormat = "@"
$wsGaz.Range("A173").Value = "2025-12-08"
$wsGaz.Range("A173").Style = "Normal"
$wsGaz.Range("B173").Value = 25.575

# ---------------------------------------------------------------------------
# 3) CO2 — append 2025-12-08 row, price not yet published.
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A173").NumberFormat = "@"
$wsCo2.Range("A173").Value = "2025-12-08"
$wsCo2.Range("A173").Style = "Normal"
$wsCo2.Range("B173").Value = ""
